# Interview 2 Talking Points.docx - apply commit edits via Word COM interop.
#
# We perform the paragraph-level restructuring bottom-to-top so the
# (1-based) $d.Paragraphs(n) indices used for the earlier edits remain
# valid while later (higher-numbered) edits are applied first.

$d = $word.ActiveDocument

function New-PkgXml([string]$bodyInner) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step A ---------------------------------------------------------
# Replace the "Authentication" / "How does one log into..." (old,
# etc etc. wording, with lastRenderedPageBreak) / "Confirm , How would
# you like to sign health form?" paragraphs with the two new closing
# paragraphs: the "Are you legally required..." question and an empty
# paragraph carrying the _GoBack bookmark.
$pAuth = $d.Paragraphs(21)
$pConfirm = $d.Paragraphs(23)
$rStart = $pAuth.Range.Start
$rEnd = $pConfirm.Range.End
$rStepA = $d.Range($rStart, $rEnd)
$stepAInner = '<w:p><w:r><w:t>Are you legally required to have a signature or is a digital signature sufficient for health forms, order delivery etc?</w:t></w:r></w:p><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$rStepA.InsertXML((New-PkgXml $stepAInner))

# --- Step B ---------------------------------------------------------
# After "...or just setup to print the qrcode?" (paragraph 19), insert
# a new (empty-text) paragraph whose mark carries single-underline
# formatting, followed by the relocated "How does one log into a
# military networked system..." question (now split into 3 runs and
# with "etc.?" instead of "etc etc.?", and no lastRenderedPageBreak).
# The following paragraph (originally an unchanged blank separator)
# must survive untouched, but InsertXML at a paragraph-boundary always
# folds the *last* fragment paragraph's mark into the paragraph that
# follows the insertion point. A trailing throwaway blank paragraph in
# the fragment absorbs that fold instead, and is then deleted again so
# the original blank separator paragraph is left standing on its own.
$p19 = $d.Paragraphs(19)
$rStepB = $d.Range($p19.Range.End, $p19.Range.End)
$stepBInner = '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p><w:p><w:r><w:t>How does one log into a military networked system, login</w:t></w:r><w:r><w:t xml:space="preserve"> via email/phone number/id/name</w:t></w:r><w:r><w:t xml:space="preserve"> etc.?</w:t></w:r></w:p><w:p/>'
$rStepB.InsertXML((New-PkgXml $stepBInner))
$pThrowaway = $d.Paragraphs(22)
$pThrowaway.Range.Delete()

# --- Step C ---------------------------------------------------------
# Drop the standalone "Printing orders or just the qr code:" paragraph.
$p18 = $d.Paragraphs(18)
$p18.Range.Delete()

# --- Step D ---------------------------------------------------------
# Insert a new blank paragraph right after "Will the application be
# installed on one device ... installed in it?" (paragraph 15).
# InsertParagraphAfter() leaves a stray empty <w:r/> behind, so we
# immediately normalise the freshly-created paragraph back down to a
# bare <w:p/> via InsertXML.
$p15 = $d.Paragraphs(15)
$p15.Range.InsertParagraphAfter()
$pNewBlank = $d.Paragraphs(16)
$stepDInner = '<w:p/><w:p/>'
$pNewBlank.Range.InsertXML((New-PkgXml $stepDInner))

# --- Step E ---------------------------------------------------------
# "(show database table) image attached" becomes an empty paragraph
# (the paragraph itself stays, only its run content is removed).
$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
$r9Text = $d.Range($r9.Start, $r9.End - 1)
$r9Text.Delete()

# --- Step F ---------------------------------------------------------
# Tag the drawing run (the picture just above) with an explicit
# Far-East language of en-AU -> <w:lang w:eastAsia="en-AU"/>.
$p8 = $d.Paragraphs(8)
$p8.Range.LanguageIDFarEast = "en-AU"
